$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("R40") in the rules table: the "From" value (column C) changes from 18 to 1.
$ws.Range("C10").Value = 1
